$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the product in row 2: old item (2845959 / Pril ISIS Cold Power liquid 650ml Lemon)
# is swapped for a new item (2970482 / Pril Isis Ultra Power 650ml).
# The new Item ID is stored as text (not a number) in the edited workbook, so force a
# leading quote prefix the way Excel itself would when a numeric-looking value is typed
# into a cell as text.
$ws.Range("A2").Value = "'2970482"
$ws.Range("B2").Value = "Pril Isis Ultra Power 650ml"

# Reflect the new selection left behind by the editor.
$ws.Range("A2:B2").Select()
